$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.728132367134094
$ws.Range("B1").Value = 3.544512510299683
$ws.Range("C1").Value = 2.966129302978516
$ws.Range("D1").Value = 1.999465703964233
$ws.Range("E1").Value = 1.162101864814758
